$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.120266318321228
$ws.Range("B1").Value = 1.705551743507385
$ws.Range("C1").Value = 6.978481292724609
$ws.Range("D1").Value = 2.712957859039307
$ws.Range("E1").Value = 1.460768103599548
